$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4.75
$ws.Range("H2").Value = 4.2
$ws.Range("I2").Value = 1.62
$ws.Range("J2").Value = 4.75
$ws.Range("L2").Value = 2.1
$ws.Range("M2").Value = 1.02
$ws.Range("N2").Value = 19
$ws.Range("Q2").Value = 1.5
$ws.Range("R2").Value = 2.5
$ws.Range("S2").Value = 1.25
$ws.Range("T2").Value = 3.75
$ws.Range("W2").Value = 19
$ws.Range("X2").Value = 29
$ws.Range("AC2").Value = 19
$ws.Range("AD2").Value = 8.5
$ws.Range("AK2").Value = 13
$ws.Range("AM2").Value = 19
$ws.Range("AN2").Value = 7
$ws.Range("AO2").Value = 23
$ws.Range("AS2").Value = 126
$ws.Range("AT2").Value = 3.75
$ws.Range("AY2").Value = 8
$ws.Range("AZ2").Value = 15
$ws.Range("BA2").Value = 23
$ws.Range("BC2").Value = 81
# Row 3
$ws.Range("O3").Value = 1.11
$ws.Range("P3").Value = 6.5
$ws.Range("Q3").Value = 1.4
$ws.Range("R3").Value = 2.88
# Row 4
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 3.5
$ws.Range("I4").Value = 3.75
$ws.Range("J4").Value = 2.63
$ws.Range("L4").Value = 4
$ws.Range("Q4").Value = 1.83
$ws.Range("R4").Value = 1.98
$ws.Range("W4").Value = 8
$ws.Range("Y4").Value = 9
$ws.Range("AF4").Value = 41
$ws.Range("AH4").Value = 12
$ws.Range("AO4").Value = 11
$ws.Range("AX4").Value = 5.5
$ws.Range("AZ4").Value = 26
$ws.Range("BC4").Value = 151
# Row 5
$ws.Range("G5").Value = 1.62
$ws.Range("H5").Value = 3.6
$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 2.3
$ws.Range("U5").Value = 2.63
$ws.Range("V5").Value = 1.44
$ws.Range("Z5").Value = 11
$ws.Range("AH5").Value = 11
$ws.Range("AI5").Value = 29
$ws.Range("AK5").Value = 81
$ws.Range("AN5").Value = 3.25
$ws.Range("AO5").Value = 9
$ws.Range("AP5").Value = 29
$ws.Range("AU5").Value = 11
$ws.Range("AV5").Value = 101
$ws.Range("AX5").Value = 7.5
$ws.Range("AY5").Value = 41
$ws.Range("AZ5").Value = 51
# Row 6
$ws.Range("G6").Value = 1.9
$ws.Range("I6").Value = 4.1
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 10
$ws.Range("Q6").Value = 1.95
$ws.Range("R6").Value = 1.9
$ws.Range("X6").Value = 9
$ws.Range("AC6").Value = 10
$ws.Range("AF6").Value = 51
$ws.Range("AL6").Value = 34
$ws.Range("AM6").Value = 41
$ws.Range("AW6").Value = 151
$ws.Range("AX6").Value = 6
$ws.Range("BD6").Value = 151
# Row 7
$ws.Range("H7").Value = 3.8
$ws.Range("I7").Value = 1.73
$ws.Range("K7").Value = 2.25
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 11
$ws.Range("Q7").Value = 1.85
$ws.Range("R7").Value = 2
$ws.Range("S7").Value = 1.36
$ws.Range("T7").Value = 3
$ws.Range("AB7").Value = 41
$ws.Range("AC7").Value = 12
$ws.Range("AD7").Value = 7
$ws.Range("AF7").Value = 51
$ws.Range("AH7").Value = 7.5
$ws.Range("AR7").Value = 101
$ws.Range("AT7").Value = 3
$ws.Range("BB7").Value = 51
# Row 8
$ws.Range("G8").Value = 3.2
$ws.Range("I8").Value = 2.2
$ws.Range("J8").Value = 3.75
$ws.Range("M8").Value = 1.06
$ws.Range("N8").Value = 10
$ws.Range("Q8").Value = 1.95
$ws.Range("R8").Value = 1.9
$ws.Range("S8").Value = 1.36
$ws.Range("T8").Value = 3
$ws.Range("W8").Value = 11
$ws.Range("AT8").Value = 3
# Row 9
$ws.Range("O9").Value = 1.4
$ws.Range("P9").Value = 3
$ws.Range("Q9").Value = 2.2
$ws.Range("R9").Value = 1.67
# Row 10
$ws.Range("G10").Value = 1.9
$ws.Range("I10").Value = 4.2
$ws.Range("J10").Value = 2.6
$ws.Range("X10").Value = 8.5
$ws.Range("Y10").Value = 9
$ws.Range("AL10").Value = 34
$ws.Range("AW10").Value = 151
$ws.Range("BB10").Value = 101
# Row 11
$ws.Range("G11").Value = 1.75
$ws.Range("H11").Value = 3.7
$ws.Range("I11").Value = 4.5
$ws.Range("J11").Value = 2.38
$ws.Range("K11").Value = 2.25
$ws.Range("L11").Value = 4.75
$ws.Range("M11").Value = 1.05
$ws.Range("N11").Value = 11
$ws.Range("Q11").Value = 1.85
$ws.Range("Z11").Value = 15
$ws.Range("AC11").Value = 11
$ws.Range("AD11").Value = 7
$ws.Range("AH11").Value = 13
$ws.Range("AI11").Value = 23
$ws.Range("AL11").Value = 34
$ws.Range("AO11").Value = 9
$ws.Range("AQ11").Value = 29
$ws.Range("AR11").Value = 51
$ws.Range("AX11").Value = 6
$ws.Range("AY11").Value = 23
# Row 14
$ws.Range("H14").Value = 3.25
$ws.Range("I14").Value = 2.72
$ws.Range("J14").Value = 2.92
$ws.Range("K14").Value = 2.18
$ws.Range("L14").Value = 3.25
$ws.Range("N14").Value = 8
$ws.Range("S14").Value = 1.34
$ws.Range("T14").Value = 3
$ws.Range("V14").Value = 2.3
$ws.Range("W14").Value = 9.75
$ws.Range("X14").Value = 13
$ws.Range("Z14").Value = 26
$ws.Range("AA14").Value = 18
$ws.Range("AB14").Value = 23
$ws.Range("AC14").Value = 8
$ws.Range("AD14").Value = 6.6
$ws.Range("AE14").Value = 11.5
$ws.Range("AH14").Value = 11.25
$ws.Range("AI14").Value = 16.5
$ws.Range("AL14").Value = 21
$ws.Range("AM14").Value = 24
$ws.Range("AP14").Value = 18
$ws.Range("AR14").Value = 70
$ws.Range("AS14").Value = 200
$ws.Range("AT14").Value = 3
$ws.Range("AU14").Value = 6.4
$ws.Range("AX14").Value = 4.9
$ws.Range("AY14").Value = 14.5
$ws.Range("AZ14").Value = 18.5
$ws.Range("BA14").Value = 60
$ws.Range("BB14").Value = 80
